# Rewrites the "estado de cuenta" worker table (rows 16-56) with the updated
# data set described in the commit: old periods/workers removed, new rows
# added for CINTHIA HERAZO CASTILLO, AUGUSTO SNEIDER SOMERSON RAMIREZ grows to
# period 2405, and the "Salario Basico" for MAYCOL SOLANO REYES drops to
# 1000000 (was 1300000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Each worker's block of 6 periods (2412 down to 2407), in this order.
$workers = @(
    @{ Doc = "1143414285"; Name = "CINTHIA HERAZO CASTILLO";           Salario = 1300000 },
    @{ Doc = "1004367229"; Name = "AUGUSTO SNEIDER SOMERSON RAMIREZ";  Salario = 1300000 },
    @{ Doc = "1001976903"; Name = "MAYCOL SOLANO REYES";               Salario = 1000000 },
    @{ Doc = "1063153542"; Name = "WILLIAM JOSE SANCHEZ GABALO";       Salario = 1300000 },
    @{ Doc = "1007375876"; Name = "DANIELA PAOLA BARRIOS MIRANDA";     Salario = 1300000 }
)

$periodsLong  = @("2412", "2411", "2410", "2409", "2408", "2407", "2406", "2405")
$periodsShort = @("2412", "2411", "2410", "2409", "2408", "2407")

$row = 16

foreach ($w in $workers) {
    if ($w.Name -eq "AUGUSTO SNEIDER SOMERSON RAMIREZ" -or $w.Name -eq "MAYCOL SOLANO REYES" -or $w.Name -eq "WILLIAM JOSE SANCHEZ GABALO") {
        $periods = $periodsLong
    } else {
        $periods = $periodsShort
    }

    foreach ($p in $periods) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $w.Doc
        $ws.Cells.Item($row, 4).Value = $w.Name
        $ws.Cells.Item($row, 5).Value = $p
        $ws.Cells.Item($row, 6).Value = 52000
        $ws.Cells.Item($row, 7).Value = $w.Salario
        $row = $row + 1
    }
}

# Final "2501" period row for every worker, in the same order.
foreach ($w in $workers) {
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = $w.Doc
    $ws.Cells.Item($row, 4).Value = $w.Name
    $ws.Cells.Item($row, 5).Value = "2501"
    $ws.Cells.Item($row, 6).Value = 50266
    $ws.Cells.Item($row, 7).Value = $w.Salario
    $row = $row + 1
}
